$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of player rows (2-19). Each entry is Name, Position, Team.
$data = @(
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Toumani Camara", "SF,PF", "Portland Trail Blazers"),
    @("Nikola Jovic", "PF,C", "Miami Heat"),
    @("Klay Thompson", "SG,SF", "Dallas Mavericks"),
    @("Spencer Dinwiddie", "PG,SG", "Dallas Mavericks"),
    @("Dennis Schröder", "PG,SG", "Golden State Warriors"),
    @("Trayce Jackson-Davis", "PF,C", "Golden State Warriors"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Joel Embiid", "C", "Philadelphia 76ers")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
